try {
    $p = $ppt.ActivePresentation

    # ------------------------------------------------------------------
    # 1) Refresh the "Date Placeholder" fields (master + every slide
    #    layout) from 4/13/25 -> 4/17/25.
    # ------------------------------------------------------------------
    function Update-DatePlaceholders($shapes, $newText) {
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.Name -like "Date Placeholder*") {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }

    $master = $p.SlideMaster
    Update-DatePlaceholders $master.Shapes "4/17/25"

    $layouts = $master.CustomLayouts
    for ($li = 1; $li -le $layouts.Count; $li++) {
        $layout = $layouts.Item($li)
        Update-DatePlaceholders $layout.Shapes "4/17/25"
    }

    # ------------------------------------------------------------------
    # 2) Slide 3, "Kotak Teks 4": "... menampilkan pengguna yang ..."
    #    becomes "... menampilkan dashboard pengguna yang ...".
    # ------------------------------------------------------------------
    $slide3 = $p.Slides.Item(3)
    $targetShape = $null
    for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
        $shp = $slide3.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "*Halaman Admin*pengguna*") {
                $targetShape = $shp
            }
        }
    }

    $tr = $targetShape.TextFrame.TextRange
    $full = $tr.Text
    $needle = " pengguna"
    $idx = $full.IndexOf($needle)
    if ($idx -ge 0) {
        $startOneBased = $idx + 1
        $chars = $tr.Characters($startOneBased, $needle.Length)
        $chars.Text = " dashboard pengguna"
    }

    Write-Host "Edit applied OK"
} catch {
    Write-Host "ERROR: $_"
}
